# Add a new "GroepClassificatieCode" column (C) with group labels next to
# the existing Classificatiecode / Groepnaam table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "GroepClassificatieCode"
$ws.Range("C2").Value = "groep1"
$ws.Range("C3").Value = "groep2"
$ws.Range("C4").Value = "groep2"
$ws.Range("C5").Value = "groep3"
$ws.Range("C6").Value = "groep4"
$ws.Range("C7").Value = "groep5"
$ws.Range("C8").Value = "groep6"
$ws.Range("C9").Value = "groep7"

# Size column C so it fits its contents (matches the 22-wide bestFit column
# the author ended up with).
$ws.Columns("C:C").ColumnWidth = 21.17

# Match the author's final cell selection when the file was saved.
[void]$ws.Range("F11").Select()
